# ---------------------------------------------------------------------
# B1--and-B2-PowerPoint.pptx edit
#
# 1) The table on slide 5 switches from table style
#    {CFDB9227-52E8-4961-BEEC-856CD7F97FB1} to the built-in style
#    {5D84E2F9-6312-4FBD-8B61-08E04E047A97} ("Medium Style 2 - Accent 1").
# 2) The deck's theme colour scheme ("Integral" / "Red Violet") is
#    replaced with the standard Office colour scheme (dk1, lt1, dk2,
#    lt2, accent1-6, hlink, folHlink).
# ---------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1) table style -----------------------------------------------
$tableSlide = $p.Slides.Item(5)
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $shp = $tableSlide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{5D84E2F9-6312-4FBD-8B61-08E04E047A97}")
    }
}

# --- 2) theme colours -----------------------------------------------
# ThemeColorScheme.Colors is 1-based, ordered:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = $officeColors[$i - 1]
}
